$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "25.818.86"
Set-TextValue $ws.Range("E2") "  -0.06%  "
Set-TextValue $ws.Range("D3") "1.640.40"
Set-TextValue $ws.Range("E3") "  +0.72%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.24%  "
Set-TextValue $ws.Range("D5") "215.62"
Set-TextValue $ws.Range("E5") "  +0.04%  "
Set-TextValue $ws.Range("D6") "0.5061"
Set-TextValue $ws.Range("E6") "  +0.19%  "
Set-TextValue $ws.Range("D7") "1.005"
Set-TextValue $ws.Range("E7") "  +0.27%  "
Set-TextValue $ws.Range("D8") "0.2586"
Set-TextValue $ws.Range("E8") "  +0.83%  "
Set-TextValue $ws.Range("D9") "0.06425"
Set-TextValue $ws.Range("E9") "  +1.49%  "
Set-TextValue $ws.Range("D10") "20.29"
Set-TextValue $ws.Range("E10") "  +4.32%  "
Set-TextValue $ws.Range("D11") "0.07814"
Set-TextValue $ws.Range("E11") "  +0.79%  "
Set-TextValue $ws.Range("D12") "4.290"
Set-TextValue $ws.Range("E12") "  +1.37%  "
Set-TextValue $ws.Range("D13") "1.871.48"
Set-TextValue $ws.Range("E13") "  +0.98%  "
Set-TextValue $ws.Range("D14") "1.637.70"
Set-TextValue $ws.Range("E14") "  +0.62%  "
Set-TextValue $ws.Range("D15") "0.5630"
Set-TextValue $ws.Range("E15") "  +2.77%  "
Set-TextValue $ws.Range("D16") "0.0₅7652"
Set-TextValue $ws.Range("E16") "  +0.17%  "
Set-TextValue $ws.Range("D17") "63.30"
Set-TextValue $ws.Range("E17") "  -0.52%  "
Set-TextValue $ws.Range("D18") "25.874.19"
Set-TextValue $ws.Range("E18") "  +0.10%  "
Set-TextValue $ws.Range("E19") "  +0.11%  "
Set-TextValue $ws.Range("D20") "193.65"
Set-TextValue $ws.Range("E20") "  +0.10%  "
Set-TextValue $ws.Range("D21") "4.359"
Set-TextValue $ws.Range("E21") "  -1.37%  "
Set-TextValue $ws.Range("D22") "9.929"
Set-TextValue $ws.Range("E22") "  +0.55%  "
Set-TextValue $ws.Range("D23") "6.104"
Set-TextValue $ws.Range("E23") "  +1.13%  "
Set-TextValue $ws.Range("D24") "1.003"
Set-TextValue $ws.Range("E24") "  +0.14%  "
Set-TextValue $ws.Range("D25") "1.799"
Set-TextValue $ws.Range("E25") "  -6.45%  "
Set-TextValue $ws.Range("D26") "140.26"
Set-TextValue $ws.Range("E26") "  -1.22%  "
Set-TextValue $ws.Range("D27") "0.1247"
Set-TextValue $ws.Range("E27") "  +1.07%  "
Set-TextValue $ws.Range("D28") "6.824"
Set-TextValue $ws.Range("E28") "  +0.65%  "
Set-TextValue $ws.Range("D29") "15.50"
Set-TextValue $ws.Range("E29") "  -0.12%  "
Set-TextValue $ws.Range("D30") "1.245"
Set-TextValue $ws.Range("E30") "  +0.52%  "
Set-TextValue $ws.Range("D31") "0.04927"
Set-TextValue $ws.Range("E31") "  +1.12%  "
Set-TextValue $ws.Range("D32") "3.311"
Set-TextValue $ws.Range("E32") "  +2.26%  "
Set-TextValue $ws.Range("D33") "3.247"
Set-TextValue $ws.Range("E33") "  +2.18%  "
Set-TextValue $ws.Range("D34") "1.584"
Set-TextValue $ws.Range("E34") "  +3.07%  "
Set-TextValue $ws.Range("D35") "2.384"
Set-TextValue $ws.Range("E35") "  +0.58%  "
Set-TextValue $ws.Range("D36") "0.9074"
Set-TextValue $ws.Range("E36") "  +1.34%  "
Set-TextValue $ws.Range("D37") "2.581"
Set-TextValue $ws.Range("E37") "  +1.26%  "
Set-TextValue $ws.Range("D38") "0.5553"
Set-TextValue $ws.Range("E38") "  +1.25%  "
Set-TextValue $ws.Range("D39") "1.128.83"
Set-TextValue $ws.Range("E39") "  +0.82%  "
Set-TextValue $ws.Range("D40") "0.01566"
Set-TextValue $ws.Range("E40") "  +1.09%  "
Set-TextValue $ws.Range("D41") "1.002"
Set-TextValue $ws.Range("D42") "5.536"
Set-TextValue $ws.Range("E42") "  -0.34%  "
Set-TextValue $ws.Range("D43") "0.8020"
Set-TextValue $ws.Range("E43") "  +0.78%  "
Set-TextValue $ws.Range("D44") "98.30"
Set-TextValue $ws.Range("E44") "  +1.34%  "
Set-TextValue $ws.Range("D45") "1.781.62"
Set-TextValue $ws.Range("E45") "  +1.07%  "
Set-TextValue $ws.Range("E46") "  -7.28%  "
Set-TextValue $ws.Range("D47") "55.64"
Set-TextValue $ws.Range("E47") "  +1.93%  "
Set-TextValue $ws.Range("E48") "  -3.55%  "
Set-TextValue $ws.Range("D49") "7.767"
Set-TextValue $ws.Range("E49") "  +3.22%  "
Set-TextValue $ws.Range("D50") "0.05034"
Set-TextValue $ws.Range("E50") "  -1.81%  "
Set-TextValue $ws.Range("D51") "1.002"
Set-TextValue $ws.Range("E51") "  -0.03%  "
